# Daily attendance processing - reorders the "Recorded By" list in column G
# so that the most-recently-recording entry (the last item in the
# comma-separated list) is moved to the front of the list.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Rows whose "Recorded By" value is intentionally left untouched by the
# daily processing run (already in the desired order / not reprocessed).
$skipRows = @(7, 33, 59)

for ($r = $firstRow; $r -le $lastRow; $r++) {

    if ($skipRows -contains $r) {
        continue
    }

    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($null -eq $raw) {
        continue
    }

    $text = [string]$raw
    if ($text.Length -eq 0) {
        continue
    }

    if ($text.IndexOf(",") -lt 0) {
        # Only a single recorder listed - nothing to rotate.
        continue
    }

    $parts = $text.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $count = $trimmed.Length
    if ($count -lt 2) {
        continue
    }

    $lastItem = $trimmed[$count - 1]
    $remaining = $trimmed[0..($count - 2)]
    $newParts = @($lastItem) + $remaining
    $newText = [string]::Join(", ", $newParts)

    if ($newText -ne $text) {
        $cell.Value2 = $newText
    }
}
